$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add a "2021" column (R) mirroring the formatting of the existing
#        "2020" column (Q) for the data rows (3-37), then fill in the new
#        2021 values cell by cell. Copy/PasteSpecial(formats) first so each
#        new R cell inherits the same style as its Q neighbour, then the
#        values are written on top (row 34 stays a blank spacer row, so it
#        only gets the format, never a value). ---
$ws.Range("Q3:Q37").Copy()
$ws.Range("R3:R37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 0.12641839647678207
$ws.Range("R5").Value = 0.14922981985616976
$ws.Range("R6").Value = 0.10326895933792253
$ws.Range("R7").Value = 0.03433011112114915
$ws.Range("R8").Value = 0.036820478077087354
$ws.Range("R9").Value = 0.031930519190242035
$ws.Range("R10").Value = 0.087302929367211068
$ws.Range("R11").Value = 0.10296328329317765
$ws.Range("R12").Value = 0.071859056271889668
$ws.Range("R13").Value = 0.10716050460690947
$ws.Range("R14").Value = 0.079035451351703812
$ws.Range("R15").Value = 0.13553052227085377
$ws.Range("R16").Value = 0.06479643687803946
$ws.Range("R17").Value = 0.07643825526207898
$ws.Range("R18").Value = 0.053576570965516782
$ws.Range("R19").Value = 0.054163459619715498
$ws.Range("R20").Value = 0.064872252119520635
$ws.Range("R21").Value = 0.043693418784505472
$ws.Range("R22").Value = 0.051373884452794741
$ws.Range("R23").Value = 0.029662368095156877
$ws.Range("R24").Value = 0.072642215296997686
$ws.Range("R25").Value = 0.13772601093442507
$ws.Range("R26").Value = 0.15668565643254884
$ws.Range("R27").Value = 0.11816042869432726
$ws.Range("R28").Value = 0.33417383115107696
$ws.Range("R29").Value = 0.41139191068108794
$ws.Range("R30").Value = 0.24697746624641295
$ws.Range("R31").Value = 0.16773611144997194
$ws.Range("R32").Value = 0.1959922553363346
$ws.Range("R33").Value = 0.13791201213625709
$ws.Range("R35").Value = 0
$ws.Range("R36").Value = 0.1
$ws.Range("R37").Value = 0.2

# --- 2. Move the selection highlight from P30 to C1, matching the saved
#        view state in the workbook. ---
$ws.Range("C1").Select()
